# Update column G ("K" - strikeouts) values for rows 2-28 per regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 7
    3  = 5
    4  = 8
    5  = 5
    6  = 4
    7  = 9
    8  = 8
    9  = 4
    10 = 10
    11 = 5
    12 = 9
    13 = 7
    14 = 8
    15 = 8
    16 = 8
    17 = 8
    18 = 6
    19 = 11
    20 = 6
    21 = 6
    22 = 3
    23 = 6
    24 = 10
    25 = 6
    26 = 5
    27 = 2
    28 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
